$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update summary header values ---
$ws.Range("E11").Value = 324700   # VALOR MORA total
$ws.Range("F13").Value = 6        # Cant. Periodos

# --- Insert a new row in the detail table, then copy formatting from the row above ---
$ws.Rows.Item(17).Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats

# --- Rewrite the detail rows (16-21) with final data ---
# Row 16: MIGUEL SALCEDO OÑATE moves to the top, with updated mora value
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "84104293"
$ws.Range("D16").Value = "MIGUEL SALCEDO OÑATE"
$ws.Range("E16").Value = "2208"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1423500

# Rows 17-21: JHURGEN ARLEYS MOLINARES PADILLA, periods 2504..2508 ascending
$periods = @(2504, 2505, 2506, 2507, 2508)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1143382239"
    $ws.Cells.Item($r, 4).Value = "JHURGEN ARLEYS MOLINARES PADILLA"
    $ws.Cells.Item($r, 5).Value = [string]$periods[$i]
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 689455
}
